$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New recruitment rows appended after the existing data (rows 2-46).
# Columns: A=Date, B=Center, C=Intervention, D=Uptake, E=E-cigarettes,
#          F=Nicotine pouches, G=Nicotine Patches, H=Control

$rows = @(
  @{ Row=47; Date=45811; Center="Basel";      C=1;    D=1;    E=$null; F=$null; G=1;    H=$null },
  @{ Row=48; Date=45811; Center="Basel";      C=1;    D=1;    E=$null; F=$null; G=1;    H=$null },
  @{ Row=49; Date=45812; Center="St.Gallen";  C=1;    D=$null;E=$null; F=$null; G=$null;H=1 },
  @{ Row=50; Date=45812; Center="Basel";      C=$null;D=$null;E=$null; F=$null; G=$null;H=1 },
  @{ Row=51; Date=45813; Center="Basel";      C=1;    D=1;    E=$null; F=$null; G=1;    H=$null },
  @{ Row=52; Date=45813; Center="Basel";      C=1;    D=1;    E=1;     F=$null; G=$null;H=$null }
)

$ws.Range("A46").Copy()
$ws.Range("A47:A52").PasteSpecial(-4122)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $r.Date
  $ws.Cells.Item($row, 2).Value = $r.Center
  if ($r.C -ne $null) { $ws.Cells.Item($row, 3).Value = $r.C }
  if ($r.D -ne $null) { $ws.Cells.Item($row, 4).Value = $r.D }
  if ($r.E -ne $null) { $ws.Cells.Item($row, 5).Value = $r.E }
  if ($r.F -ne $null) { $ws.Cells.Item($row, 6).Value = $r.F }
  if ($r.G -ne $null) { $ws.Cells.Item($row, 7).Value = $r.G }
  if ($r.H -ne $null) { $ws.Cells.Item($row, 8).Value = $r.H }
}

$ws.Range("F46").Select()
